$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(28.51541462446781,25.51142295962591,31.7991150398072),
    @(22.25585576105784,15.97302237025835,28.66912699016356),
    @(23.05137957723008,17.82749041943472,28.02240568169178),
    @(27.05833038679527,20.08663067962028,34.2828205169893),
    @(30.36072101474461,26.44783789281137,34.59189900390377),
    @(24.37532435074546,21.60866166128789,27.18162748542422),
    @(18.57787027929137,16.32046993974671,20.58476732696001),
    @(28.11279401053422,24.09200286578591,32.07790012100787),
    @(20.16732719625286,17.95955450840523,22.40156183391565),
    @(27.18587506358066,24.29430315253474,30.43470810300648),
    @(34.08247422477606,30.27944907125182,37.79429511272),
    @(25.79867440168809,22.61251849102035,29.08675109757067),
    @(25.78301001085314,22.8566237706098,28.97182714690977),
    @(31.90421718790628,23.50875671478268,39.77631477124827),
    @(30.34625591636008,26.46846294177395,34.63785526380993),
    @(34.05810921320352,30.4392938455772,37.69952022405466),
    @(22.99588535120866,17.63798837506834,27.99345602621706),
    @(30.80462091641987,27.01525597411609,34.68717226359928),
    @(34.9473738589548,30.38752113162489,39.60989701792482),
    @(32.75506573647809,27.11668045086042,37.32843546979694),
    @(21.80995717900162,18.93776961292972,24.64380067894182),
    @(34.82257063613874,27.5341585525012,42.07012057842832),
    @(30.39168735149386,26.51715560414863,34.68725751168721),
    @(25.61120525301443,22.09335832702901,29.70619369243233),
    @(32.71103480598697,26.46723420982992,37.97485425997139),
    @(28.67089261496256,25.56840679581229,32.03518736716897),
    @(26.46687990905067,23.52109760193689,29.42048037447858),
    @(36.72459776692349,29.08878502376017,43.7303415509771),
    @(31.92652740531376,28.49314262722579,35.23039431621249),
    @(23.98535853781496,20.79367476935354,26.91687765369139),
    @(21.08124153907659,18.662561082002,23.65620670138344),
    @(25.67532693977791,22.64679942108379,28.86830485503478),
    @(34.36223268485242,30.70663890948799,38.50026941071764),
    @(32.57906889541039,28.45223320945659,36.91773690777674),
    @(24.41196351058469,21.60834160200385,27.23116296881804),
    @(35.41118491982756,29.68441667383447,40.93687145152457),
    @(28.83161577189314,25.66492001892161,32.35612546067367),
    @(35.8555231823269,31.87649645875045,39.72485114335048),
    @(27.987184993619,23.07361564754637,32.86890422776501),
    @(22.42246372837064,19.89414530412046,25.20593907512853),
    @(28.86433194979812,25.67399500392279,32.14947489185039),
    @(31.90149893669372,28.42758337975143,35.22510106079826),
    @(21.36480862600348,18.66260073062385,24.18323627639434),
    @(32.03568875626988,28.07543052794952,35.86744093925697),
    @(21.97844619102932,19.60179687370963,24.62062785779077),
    @(23.45557266864113,20.75414246447868,26.54507026952158),
    @(34.04304821481757,30.20978619686585,37.75068216899476),
    @(25.74833846352197,21.03496185819164,30.92795723557029),
    @(35.38566221054234,29.63575608392734,40.89417186443777),
    @(22.28272382817056,16.02596892543221,28.67372154748251),
    @(20.38592705416354,17.84167645736595,23.34262459278627),
    @(35.92000503397135,32.01397488309934,39.63249876329957),
    @(31.90714649542242,23.53193615658738,39.80576543179825),
    @(32.78331516705319,27.6695494013407,38.03219964611328),
    @(23.9439608955713,21.31372376851823,26.67906419719397),
    @(21.07896347036549,18.67951945920808,23.64321572038936),
    @(35.44095702616701,29.76550065509952,40.69011109332545),
    @(25.11558039467248,22.17321500147788,28.44651881624265),
    @(34.9220590779423,27.67000236912793,42.45642769960442),
    @(27.18941033937773,24.05069727602547,30.15998240707067),
    @(27.12062629274434,24.24139091299088,30.36898203933858),
    @(30.29709564388545,26.68156986817262,34.01539777111161),
    @(34.93643977085534,27.66467961789436,42.3527481044765),
    @(19.2879347096092,16.91737711187096,21.74451204522518),
    @(32.45288066287506,28.58991557998448,36.30700802771688),
    @(30.39865854329875,27.04289746513122,33.57940980553768),
    @(30.48126474779218,25.47849874643118,35.66315275719328),
    @(23.33309178856802,20.46892466202684,26.2616795960842),
    @(28.19854607549327,22.77045372073315,33.81401160640677),
    @(22.88830681871856,19.89358243355528,26.13214250315592),
    @(27.58274790002504,24.37758645660139,30.98902035058762)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
